# Apply the crypto-price refresh for Sheet1 (rows 2-51).
# Row 34 gains a new "Frax" entry; every row from 34 downward
# shifts by one, and the previous last row ("Aave") drops off.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text for rows whose new value would
# otherwise be auto-parsed as a number (losing trailing zeros, etc).
$ws.Range("D4:D9").NumberFormat = "@"
$ws.Range("D11:D14").NumberFormat = "@"
$ws.Range("D19:D20").NumberFormat = "@"
$ws.Range("D22:D25").NumberFormat = "@"
$ws.Range("D27:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.426.92"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.726.89"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "242.99"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4882"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").Value = "0.2618"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").Value = "0.06201"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "1.729.18"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "0.07017"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "15.48"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "4.561"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "0.5985"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "77.41"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "26.439.88"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "0.000007162"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").Value = "11.43"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").Value = "1.950.87"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "4.489"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "8.599"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("D24").Value = "5.183"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "138.95"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("D26").Value = "15.28"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").Value = "1.407"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "106.97"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "1.724"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").Value = "3.965"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "0.07988"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "3.684"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "0.04530"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "0.9997"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.614"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.6234"
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "0.9081"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "1.982"
$ws.Range("E39").Value = "  -5.17%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.407"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.01491"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "100.39"
$ws.Range("E43").Value = "  -4.18%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.411"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3872"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.700"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1156"
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05364"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.31"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.728"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.255"
$ws.Range("E51").Value = "  -0.88%  "
